$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Processes sheet: rotate columns B/C/D so the table reads
#    key | type | fuel | product   (was key | fuel | product | type)
# ---------------------------------------------------------------------------
$wsProcesses = $wb.Worksheets.Item("Processes")

for ($r = 1; $r -le 11; $r++) {
    $oldB = $wsProcesses.Cells.Item($r, 2).Value2
    $oldC = $wsProcesses.Cells.Item($r, 3).Value2
    $oldD = $wsProcesses.Cells.Item($r, 4).Value2

    $wsProcesses.Cells.Item($r, 2).Value = $oldD
    $wsProcesses.Cells.Item($r, 3).Value = $oldB
    $wsProcesses.Cells.Item($r, 4).Value = $oldC
}

# Data validation list used to live on the (old) type column D2:D11; it now
# belongs to the (new) type column B2:B11.
$wsProcesses.Range("D2:D11").Validation.Delete()
$wsProcesses.Range("B2:B11").Validation.Add(3, 1, 1, "Validate!`$B`$2:`$B`$3")

# Column D (now "product") is a bit wider than it used to be.
$wsProcesses.Columns.Item(4).ColumnWidth = 14.736979166666666

# Selection moves from the single cell C6 to the whole table B1:D11.
$wsProcesses.Range("B1:D11").Select()

# ---------------------------------------------------------------------------
# 2. Workbook-level defined name: cgam_processes now only spans columns A:B.
# ---------------------------------------------------------------------------
$wb.Names.Item("Processes!cgam_processes").RefersTo = "=Processes!`$A`$1:`$B`$5"

# ---------------------------------------------------------------------------
# 3. Rename WasteDefinition -> WasteDefinitionx and make it the active sheet
#    (it was Exergy before).
# ---------------------------------------------------------------------------
$wsWaste = $wb.Worksheets.Item("WasteDefinition")
$wsWaste.Name = "WasteDefinitionx"
$wsWaste.Activate()
